# Auto-generated edit script applying the Chocobo_Profits data refresh diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 3290662.8
$ws.Range("J19").Value = 2611
$ws.Range("L19").Value = 2611
$ws.Range("N19").Value = -2961
$ws.Range("H69").Value = 2000
$ws.Range("I69").Value = 2000
$ws.Range("K69").Value = 6000
$ws.Range("M69").Value = -5126
$ws.Range("H72").Value = 2000
$ws.Range("I72").Value = 2000
$ws.Range("K72").Value = 18000
$ws.Range("M72").Value = -13632
$ws.Range("H80").Value = 600.5161000000001
$ws.Range("I80").Value = 264.4375
$ws.Range("J80").Value = 959
$ws.Range("K80").Value = 793.3125
$ws.Range("L80").Value = 2877
$ws.Range("M80").Value = 204.6875
$ws.Range("N80").Value = -4873
$ws.Range("H83").Value = 600.5161000000001
$ws.Range("I83").Value = 264.4375
$ws.Range("J83").Value = 959
$ws.Range("K83").Value = 2379.9375
$ws.Range("L83").Value = 8631
$ws.Range("M83").Value = 2612.0625
$ws.Range("N83").Value = -18615
$ws.Range("H98").Value = 2234.5557
$ws.Range("I98").Value = 1077.125
$ws.Range("J98").Value = 3918.0908
$ws.Range("K98").Value = 1077.125
$ws.Range("L98").Value = 3918.0908
$ws.Range("M98").Value = 420.875
$ws.Range("N98").Value = -6914.0908
$ws.Range("H113").Value = 5139.615
$ws.Range("I113").Value = 2004
$ws.Range("J113").Value = 5400.9165
$ws.Range("K113").Value = 2004
$ws.Range("L113").Value = 5400.9165
$ws.Range("M113").Value = 1250
$ws.Range("N113").Value = -11908.9165
$ws.Range("H122").Value = 2234.5557
$ws.Range("I122").Value = 1077.125
$ws.Range("J122").Value = 3918.0908
$ws.Range("K122").Value = 3231.375
$ws.Range("L122").Value = 11754.2724
$ws.Range("M122").Value = -781.375
$ws.Range("N122").Value = -16654.2724
$ws.Range("H123").Value = 41520
$ws.Range("J123").Value = 41520
$ws.Range("L123").Value = 41520
$ws.Range("N123").Value = -51320
$ws.Range("H129").Value = 824.4536000000001
$ws.Range("I129").Value = 338.15384
$ws.Range("J129").Value = 899.7143
$ws.Range("K129").Value = 1014.46152
$ws.Range("L129").Value = 2699.1429
$ws.Range("M129").Value = 3985.53848
$ws.Range("N129").Value = -12699.1429
$ws.Range("H132").Value = 5125.931
$ws.Range("I132").Value = 7299.4707
$ws.Range("K132").Value = 21898.4121
$ws.Range("M132").Value = -19368.4121
$ws.Range("H137").Value = 4053.2727
$ws.Range("I137").Value = 2958
$ws.Range("K137").Value = 8874
$ws.Range("M137").Value = -6324
$ws.Range("H138").Value = 4636.37
$ws.Range("J138").Value = 6002.2163
$ws.Range("L138").Value = 18006.6489
$ws.Range("N138").Value = -28286.6489

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 9820.546
$ws.Range("I74").Value = 75000
$ws.Range("J74").Value = 3302.6
$ws.Range("K74").Value = 75000
$ws.Range("L74").Value = 3302.6
$ws.Range("M74").Value = -74126
$ws.Range("N74").Value = -5050.6
$ws.Range("H77").Value = 9820.546
$ws.Range("I77").Value = 75000
$ws.Range("J77").Value = 3302.6
$ws.Range("K77").Value = 375000
$ws.Range("L77").Value = 16513
$ws.Range("M77").Value = -370632
$ws.Range("N77").Value = -25249
$ws.Range("H88").Value = 7411929.5
$ws.Range("I88").Value = 11114861
$ws.Range("K88").Value = 11114861
$ws.Range("M88").Value = -11114455
$ws.Range("H91").Value = 7411929.5
$ws.Range("I91").Value = 11114861
$ws.Range("K91").Value = 11114861
$ws.Range("M91").Value = -11113457
$ws.Range("H122").Value = 6670.6665
$ws.Range("I122").Value = 2012
$ws.Range("J122").Value = 9000
$ws.Range("K122").Value = 6036
$ws.Range("L122").Value = 27000
$ws.Range("M122").Value = -3586
$ws.Range("N122").Value = -31900

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2029.2307
$ws.Range("I86").Value = 1761
$ws.Range("J86").Value = 2297.4614
$ws.Range("K86").Value = 1761
$ws.Range("L86").Value = 2297.4614
$ws.Range("M86").Value = -638
$ws.Range("N86").Value = -4543.4614
$ws.Range("H89").Value = 2029.2307
$ws.Range("I89").Value = 1761
$ws.Range("J89").Value = 2297.4614
$ws.Range("K89").Value = 8805
$ws.Range("L89").Value = 11487.307
$ws.Range("M89").Value = -3189
$ws.Range("N89").Value = -22719.307

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3361.4285
$ws.Range("I31").Value = 1132.1666
$ws.Range("J31").Value = 6333.778
$ws.Range("K31").Value = 1132.1666
$ws.Range("L31").Value = 6333.778
$ws.Range("M31").Value = -837.1666
$ws.Range("N31").Value = -6923.778
$ws.Range("H34").Value = 3361.4285
$ws.Range("I34").Value = 1132.1666
$ws.Range("J34").Value = 6333.778
$ws.Range("K34").Value = 1132.1666
$ws.Range("L34").Value = 6333.778
$ws.Range("M34").Value = -930.1666
$ws.Range("N34").Value = -6737.778
$ws.Range("H68").Value = 43835
$ws.Range("J68").Value = 43835
$ws.Range("L68").Value = 43835
$ws.Range("N68").Value = -45333
$ws.Range("H71").Value = 43835
$ws.Range("J71").Value = 43835
$ws.Range("L71").Value = 131505
$ws.Range("N71").Value = -138993
$ws.Range("H99").Value = 4137.75
$ws.Range("I99").Value = 2240
$ws.Range("J99").Value = 5000.364
$ws.Range("K99").Value = 2240
$ws.Range("L99").Value = 5000.364
$ws.Range("M99").Value = -742
$ws.Range("N99").Value = -7996.364
$ws.Range("H106").Value = 35000
$ws.Range("J106").Value = 35000
$ws.Range("L106").Value = 35000
$ws.Range("N106").Value = -37524
$ws.Range("H122").Value = 2355.5334
$ws.Range("I122").Value = 1220
$ws.Range("K122").Value = 3660
$ws.Range("M122").Value = -1210
$ws.Range("H126").Value = 4137.75
$ws.Range("I126").Value = 2240
$ws.Range("J126").Value = 5000.364
$ws.Range("K126").Value = 6720
$ws.Range("L126").Value = 15001.092
$ws.Range("M126").Value = -4250
$ws.Range("N126").Value = -19941.092
$ws.Range("H132").Value = 1789.4773
$ws.Range("I132").Value = 1269.4474
$ws.Range("J132").Value = 5083
$ws.Range("K132").Value = 3808.3422
$ws.Range("L132").Value = 15249
$ws.Range("M132").Value = -1278.3422
$ws.Range("N132").Value = -20309

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1196.6154
$ws.Range("I5").Value = 349.75
$ws.Range("J5").Value = 1922.5
$ws.Range("K5").Value = 1049.25
$ws.Range("L5").Value = 5767.5
$ws.Range("M5").Value = -937.25
$ws.Range("N5").Value = -5991.5
$ws.Range("H63").Value = 4739.8667
$ws.Range("I63").Value = 4583
$ws.Range("J63").Value = 4844.4443
$ws.Range("K63").Value = 13749
$ws.Range("L63").Value = 14533.3329
$ws.Range("M63").Value = -13000
$ws.Range("N63").Value = -16031.3329
$ws.Range("H66").Value = 4739.8667
$ws.Range("I66").Value = 4583
$ws.Range("J66").Value = 4844.4443
$ws.Range("K66").Value = 41247
$ws.Range("L66").Value = 43599.9987
$ws.Range("M66").Value = -37503
$ws.Range("N66").Value = -51087.9987
$ws.Range("H113").Value = 629.3606600000001
$ws.Range("I113").Value = 621.6957
$ws.Range("J113").Value = 652.86664
$ws.Range("K113").Value = 1865.0871
$ws.Range("L113").Value = 1958.59992
$ws.Range("M113").Value = 304.9129
$ws.Range("N113").Value = -6298.59992
$ws.Range("H131").Value = 6024910.5
$ws.Range("I131").Value = 100000400
$ws.Range("J131").Value = 840.96155
$ws.Range("K131").Value = 300001200
$ws.Range("L131").Value = 2522.88465
$ws.Range("M131").Value = -299996160
$ws.Range("N131").Value = -12602.88465
$ws.Range("H135").Value = 1196.6154
$ws.Range("I135").Value = 349.75
$ws.Range("J135").Value = 1922.5
$ws.Range("K135").Value = 3147.75
$ws.Range("L135").Value = 17302.5
$ws.Range("M135").Value = -612.75
$ws.Range("N135").Value = -22372.5
$ws.Range("H138").Value = 2443
$ws.Range("I138").Value = 1794.875
$ws.Range("J138").Value = 3480
$ws.Range("K138").Value = 5384.625
$ws.Range("L138").Value = 10440
$ws.Range("M138").Value = -244.625
$ws.Range("N138").Value = -20720
$ws.Range("H140").Value = 32263.412
$ws.Range("I140").Value = 34204.875
$ws.Range("J140").Value = 1200
$ws.Range("K140").Value = 102614.625
$ws.Range("L140").Value = 3600
$ws.Range("M140").Value = -97434.625
$ws.Range("N140").Value = -13960

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2357.5293
$ws.Range("I102").Value = 1628.3077
$ws.Range("J102").Value = 4727.5
$ws.Range("K102").Value = 1628.3077
$ws.Range("L102").Value = 4727.5
$ws.Range("M102").Value = -6.307700000000068
$ws.Range("N102").Value = -7971.5
$ws.Range("H122").Value = 2819.6
$ws.Range("I122").Value = 1962.3334
$ws.Range("J122").Value = 4105.5
$ws.Range("K122").Value = 5887.0002
$ws.Range("L122").Value = 12316.5
$ws.Range("M122").Value = -3437.0002
$ws.Range("N122").Value = -17216.5
$ws.Range("H139").Value = 94250
$ws.Range("J139").Value = 94250
$ws.Range("L139").Value = 94250
$ws.Range("N139").Value = -104530

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H114").Value = 39733.332
$ws.Range("J114").Value = 39733.332
$ws.Range("L114").Value = 39733.332
$ws.Range("N114").Value = -48411.332
$ws.Range("H115").Value = 38610
$ws.Range("J115").Value = 38610
$ws.Range("L115").Value = 38610
$ws.Range("N115").Value = -40960

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 18000
$ws.Range("I62").Value = 5666.6665
$ws.Range("J62").Value = 55000
$ws.Range("K62").Value = 5666.6665
$ws.Range("L62").Value = 55000
$ws.Range("M62").Value = -5042.6665
$ws.Range("N62").Value = -56248
$ws.Range("H65").Value = 18000
$ws.Range("I65").Value = 5666.6665
$ws.Range("J65").Value = 55000
$ws.Range("K65").Value = 28333.3325
$ws.Range("L65").Value = 275000
$ws.Range("M65").Value = -25213.3325
$ws.Range("N65").Value = -281240
$ws.Range("H81").Value = 1928.5714
$ws.Range("I81").Value = 1949.6666
$ws.Range("J81").Value = 1802
$ws.Range("K81").Value = 3899.3332
$ws.Range("L81").Value = 3604
$ws.Range("M81").Value = -2838.3332
$ws.Range("N81").Value = -5726
$ws.Range("H84").Value = 1928.5714
$ws.Range("I84").Value = 1949.6666
$ws.Range("J84").Value = 1802
$ws.Range("K84").Value = 19496.666
$ws.Range("L84").Value = 18020
$ws.Range("M84").Value = -14192.666
$ws.Range("N84").Value = -28628
$ws.Range("H136").Value = 5223.684
$ws.Range("I136").Value = 2788.3845
$ws.Range("J136").Value = 10500.167
$ws.Range("K136").Value = 8365.1535
$ws.Range("L136").Value = 31500.501
$ws.Range("M136").Value = -5815.1535
$ws.Range("N136").Value = -36600.501
$ws.Range("H139").Value = 37265.277
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 37265.277
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 37265.277
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -47545.277

Write-Output "Applied all changes"